$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "Homework 3"
$ws.Range("H3").Formula = "=14/13"
$ws.Range("H4").Formula = "=8/13"
$ws.Range("H6").Formula = "=14/13"
$ws.Range("H7").Formula = "=13/13"
$ws.Range("H8").Formula = "=13/13"
$ws.Range("H9").Formula = "=12/13"
$ws.Range("H11").Formula = "=13/13"
$ws.Range("H12").Formula = "=12.5/13"
$ws.Range("H13").Formula = "=13/13"
$ws.Range("H14").Formula = "=13/13"
$ws.Range("H15").Formula = "=12.5/13"
$ws.Range("H16").Formula = "=13/13"
$ws.Range("H17").Formula = "=13/13"
$ws.Range("H18").Formula = "=13/13"
$ws.Range("H19").Formula = "=13/13"

$ws.Columns.Item(8).ColumnWidth = 11.92

$ws.Range("H17").Select()
